$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Sheet bookkeeping
#    - the current "总计" sheet (sheetId 6 / sheet6.xml) becomes the
#      new "2022-Q1" detail sheet
#    - a brand-new "总计" sheet is appended right after it (gets the
#      next sheetId / a new sheet7.xml) and receives the refreshed
#      roll-up table
# ------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

$total = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q1)
$total.Name = "总计"

$fmtSrc = $wb.Worksheets.Item("2021-Q4")       # donor for the bold/bordered look

# ------------------------------------------------------------------
# 2. Rebuild "2022-Q1" as a fund-holdings detail sheet (A1:H14)
# ------------------------------------------------------------------
$q1.Range("A1:D6").ClearContents()

# bring over the bold/centered/bordered look used by every other
# quarter sheet (header row + the numbered A column)
$fmtSrc.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$fmtSrc.Range("A2:A14").Copy()
$q1.Range("A2:A14").PasteSpecial(-4122)

# header
$q1.Cells.Item(1,2).Value = "基金代码"
$q1.Cells.Item(1,3).Value = "基金名称"
$q1.Cells.Item(1,4).Value = "基金规模"
$q1.Cells.Item(1,5).Value = "股票总仓位"
$q1.Cells.Item(1,6).Value = "仓位占比"
$q1.Cells.Item(1,7).Value = "持有市值(亿元)"
$q1.Cells.Item(1,8).Value = "仓位排名"

# column B (fund codes) and D,E,F,G hold text (e.g. "011069",
# "45.07"); force text storage so Excel doesn't silently coerce them
# to numbers and strip leading zeros
$q1.Range("B2:B14").NumberFormat = "@"
$q1.Range("D2:G14").NumberFormat = "@"

$rows = @(
    @(0,  "481001", "工银核心价值混合A",               "45.07", "87.81", "4.00", "1.8028", 5),
    @(1,  "011069", "工银瑞信成长精选混合A",             "16.32", "72.98", "3.02", "0.4929", 9),
    @(2,  "001008", "工银国企改革主题股票",               "8.66",  "91.98", "3.52", "0.3048", 5),
    @(3,  "920002", "中金精选股票A",                   "3.40",  "86.43", "5.06", "0.1720", 8),
    @(4,  "010460", "兴业研究精选混合",                  "3.41",  "89.54", "3.39", "0.1156", 7),
    @(5,  "004818", "国寿安保目标策略灵活配置混合A",        "4.06",  "36.45", "2.46", "0.0999", 4),
    @(6,  "001672", "国寿安保智慧生活股票",               "3.56",  "85.91", "2.72", "0.0968", 10),
    @(7,  "011070", "工银瑞信成长精选混合C",             "1.91",  "72.98", "3.02", "0.0577", 9),
    @(8,  "004194", "招商中证1000指数增强A",            "1.76",  "94.40", "1.07", "0.0188", 7),
    @(9,  "004819", "国寿安保目标策略灵活配置混合C",        "0.57",  "36.45", "2.46", "0.0140", 4),
    @(10, "004195", "招商中证1000指数增强C",            "0.68",  "94.40", "1.07", "0.0073", 7),
    @(11, "920922", "中金精选股票C",                   "0.14",  "86.43", "5.06", "0.0071", 8),
    @(12, "960010", "工银核心价值混合H",                "0.00",  "87.81", "4.00", "__ZERO__", 5)
)

$r = 2
foreach ($row in $rows) {
    $q1.Cells.Item($r,1).Value = $row[0]
    $q1.Cells.Item($r,2).Value = $row[1]
    $q1.Cells.Item($r,3).Value = $row[2]
    $q1.Cells.Item($r,4).Value = $row[3]
    $q1.Cells.Item($r,5).Value = $row[4]
    $q1.Cells.Item($r,6).Value = $row[5]
    if ($row[6] -eq "__ZERO__") {
        $q1.Cells.Item($r,7).NumberFormat = "General"
        $q1.Cells.Item($r,7).Value = 0
    } else {
        $q1.Cells.Item($r,7).Value = $row[6]
    }
    $q1.Cells.Item($r,8).Value = $row[7]
    $r = $r + 1
}

# ------------------------------------------------------------------
# 3. Rebuild "总计" as the quarter-over-quarter roll-up (A1:D7),
#    adding the new 2022-Q1 row on top and keeping the rest as-is
# ------------------------------------------------------------------
$fmtSrc.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$fmtSrc.Range("A2:A7").Copy()
$total.Range("A2:A7").PasteSpecial(-4122)

$total.Cells.Item(1,2).Value = "日期"
$total.Cells.Item(1,3).Value = "持有数量(只)"
$total.Cells.Item(1,4).Value = "持有市值(亿元)"

$totalRows = @(
    @(0, "2022-Q1", 13, 3.19),
    @(1, "2021-Q4", 16, 4.97),
    @(2, "2021-Q3", 33, 7.93),
    @(3, "2021-Q2", 18, 7.72),
    @(4, "2021-Q1", 15, 3.53),
    @(5, "2020-Q4", 6,  1.76)
)

$r = 2
foreach ($row in $totalRows) {
    $total.Cells.Item($r,1).Value = $row[0]
    $total.Cells.Item($r,2).Value = $row[1]
    $total.Cells.Item($r,3).Value = $row[2]
    $total.Cells.Item($r,4).Value = $row[3]
    $r = $r + 1
}
